# Common: Very first support for gallery
# Append new translation rows (760-768) to the translations sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TranslationCell {
    param($Row, $Col, $Value)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.Value = $Value
    $cell.WrapText = $true
    $cell.Font.Size = 10
}

# Language column (A) is always "cs" for every new row.
foreach ($r in 760..768) {
    Set-TranslationCell $r 1 "cs"
}

# Write the label (B) / translation (C) cells in the exact order the
# shared-string table in the target workbook introduces each new string,
# so brand-new strings land at the expected indices and reused strings
# (already present in the table) resolve back to their existing entries.

Set-TranslationCell 760 2 "lab.build.image.cancel.upload"
Set-TranslationCell 760 3 "Zrušit nahrávání"

Set-TranslationCell 761 2 "lab.liquid.upload.tab"

Set-TranslationCell 763 3 "Nejsou dostupné žádné obrázky."

Set-TranslationCell 762 2 "lab.liquid.images.tab"

Set-TranslationCell 763 2 "common.gallery.no-images"

Set-TranslationCell 764 2 "lab.liquid.image.upload"
Set-TranslationCell 764 3 "Nahrajte obrázek liquidu"

Set-TranslationCell 765 2 "lab.liquid.image.upload.hint"
Set-TranslationCell 765 3 "Obrázek pomůže snadno identifikovat, o jaký liquid se jedná."

Set-TranslationCell 766 2 "lab.liquid.image.upload.started"
Set-TranslationCell 766 3 "Nahrávání bylo zahájeno."

Set-TranslationCell 767 3 "Obrázek byl úspěšně nahrán."

Set-TranslationCell 768 2 "lab.liquid.image.cancel.upload"

Set-TranslationCell 767 2 "lab.liquid.image.upload.success"

# Remaining cells that reuse already-existing shared strings.
Set-TranslationCell 761 3 "Nahrát obrázek"
Set-TranslationCell 762 3 "Obrázky"
Set-TranslationCell 768 3 "Zrušit nahrávání"

# Match the view state recorded in the saved workbook.
$ws.Cells.Item(764, 2).Select() | Out-Null
